# Seção 3 (Experimentos e Resultados)
# Row 33 used to be "Linha21 - Linha22" (B25-B26 ... M25-M26).
# It becomes "Número de anotações L26 + L27" (B26+B27 ... M26+M27),
# and the old "Verifications" label that used to sit in A32 is removed
# (A32 becomes an empty, but still styled, cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A32 no longer holds the "Verifications" caption - clear it but keep formatting.
$ws.Range("A32").ClearContents()

# A33 gets the new caption.
$ws.Range("A33").Value = "Número de anotações L26 + L27"

# Row 33 formulas now add rows 26 and 27 together (instead of subtracting row26 from row25).
$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M")
foreach ($col in $cols) {
    $ws.Range($col + "33").Formula = "=" + $col + "26+" + $col + "27"
}

# The previous formulas could go negative, so H33/I33/M33 used a "negative" (red) style.
# Sums of two non-negative rows are always >= 0, so those cells go back to the plain style.
$ws.Range("H33").Style = "Normal"
$ws.Range("I33").Style = "Normal"
$ws.Range("M33").Style = "Normal"

# The sheet view now has A33 selected/active instead of A16.
[void]$ws.Range("A33").Select()
